$d = $word.ActiveDocument

$replacements = @(
    @("880×3=", "267×4="),
    @("831×8=", "710×5="),
    @("687×2=", "118×6="),
    @("621×6=", "935×9="),
    @("474×9=", "667×9="),
    @("834×4=", "762×2="),
    @("707×2=", "178×4="),
    @("414×8=", "690×7="),
    @("861×6=", "935×2="),
    @("950×5=", "923×4="),
    @("408×5=", "907×8="),
    @("628×3=", "689×3="),
    @("801×2=", "371×7="),
    @("755×7=", "131×8="),
    @("575×3=", "191×8="),
    @("797×5=", "854×7="),
    @("533×6=", "886×6="),
    @("750×8=", "263×5="),
    @("394×3=", "363×9="),
    @("851×8=", "623×7="),
    @("317×4=", "513×7="),
    @("983×9=", "975×6="),
    @("218×6=", "136×5="),
    @("261×3=", "512×6="),
    @("300×4=", "738×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
